# Week 15 simulations added to the 49ers 2021 Team Data workbook.
# Appends newly-simulated values to the long space-separated run sequences
# (YDS / ST sheets) and bumps the week-over-week aggregate totals
# (OFF / DEF / ST / TURNS / PEN sheets) to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# YDS sheet: append week 15 per-drive/play simulation samples
# ---------------------------------------------------------------
$ws = $wb.Sheets.Item("YDS")

$ws.Range("B2").Value2 = $ws.Range("B2").Value2 + " -1 0 1 4 4 9 12 2 27 1 7 2 3 -3 2 5 3 3 1 8 6 2 2"
$ws.Range("C2").Value2 = $ws.Range("C2").Value2 + " -2 8 3 3 3 2 4 6 4 0 2 7 1 6 1 2 4 6 2 0 17 1 0 4 3"
$ws.Range("B3").Value2 = $ws.Range("B3").Value2 + " 7 5 4 2 10 21 6 1 4 3 22 20 16 10 13 8 9 5 17 19 8 25 21 3 9 12"
$ws.Range("C3").Value2 = $ws.Range("C3").Value2 + " 14 20 0 10 11 16 9 12 -4 11 12 19 3 5 27 4 13 5 12 22 29 32 26 23"

# ---------------------------------------------------------------
# OFF sheet: updated season totals
# ---------------------------------------------------------------
$ws = $wb.Sheets.Item("OFF")

$ws.Range("C2").Value2 = 182
$ws.Range("E2").Value2 = 8
$ws.Range("F2").Value2 = 50
$ws.Range("G2").Value2 = 52
$ws.Range("J2").Value2 = 23
$ws.Range("N2").Value2 = 11

$ws.Range("C3").Value2 = 122
$ws.Range("E3").Value2 = 29
$ws.Range("F3").Value2 = 74
$ws.Range("G3").Value2 = 21
$ws.Range("H3").Value2 = 27
$ws.Range("I3").Value2 = 44
$ws.Range("J3").Value2 = 31
$ws.Range("L3").Value2 = 203
$ws.Range("M3").Value2 = 134
$ws.Range("Q3").Value2 = 422

# ---------------------------------------------------------------
# DEF sheet: updated season totals
# ---------------------------------------------------------------
$ws = $wb.Sheets.Item("DEF")

$ws.Range("C2").Value2 = 156
$ws.Range("F2").Value2 = 49
$ws.Range("G2").Value2 = 41
$ws.Range("J2").Value2 = 20
$ws.Range("N2").Value2 = 23
$ws.Range("O2").Value2 = 20
$ws.Range("P2").Value2 = 9

$ws.Range("C3").Value2 = 144
$ws.Range("E3").Value2 = 28
$ws.Range("F3").Value2 = 78
$ws.Range("G3").Value2 = 26
$ws.Range("H3").Value2 = 25
$ws.Range("I3").Value2 = 45
$ws.Range("J3").Value2 = 47
$ws.Range("L3").Value2 = 235
$ws.Range("M3").Value2 = 162
$ws.Range("Q3").Value2 = 440

# ---------------------------------------------------------------
# ST sheet: updated season totals + appended run sequences
# ---------------------------------------------------------------
$ws = $wb.Sheets.Item("ST")

$ws.Range("B2").Value2 = 67
$ws.Range("D2").Value2 = 41
$ws.Range("F2").Value2 = 576
$ws.Range("G2").Value2 = 561
$ws.Range("J2").Value2 = 285
$ws.Range("K2").Value2 = 270
$ws.Range("L2").Value2 = 160

$ws.Range("B3").Value2 = 30

$ws.Range("D3").Value2 = $ws.Range("D3").Value2 + " 42 43 29 39 50"
$ws.Range("B4").Value2 = $ws.Range("B4").Value2 + " 56 65"
$ws.Range("D4").Value2 = $ws.Range("D4").Value2 + " 12 8 0 1 0"
$ws.Range("B5").Value2 = $ws.Range("B5").Value2 + " 18 18"
$ws.Range("D5").Value2 = "0 0 0 4 0 0"
$ws.Range("B6").Value2 = $ws.Range("B6").Value2 + " 16 24 23"

# ---------------------------------------------------------------
# TURNS sheet: updated season totals
# ---------------------------------------------------------------
$ws = $wb.Sheets.Item("TURNS")

$ws.Range("D3").Value2 = 9
$ws.Range("E3").Value2 = 12

# ---------------------------------------------------------------
# PEN sheet: updated season totals
# ---------------------------------------------------------------
$ws = $wb.Sheets.Item("PEN")

$ws.Range("B2").Value2 = 14
$ws.Range("B3").Value2 = 15
